$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 95, pushing the existing row 95 (and below)
# down to row 96. This duplicates row 95's formatting/content into row 96,
# which already holds the values we want to keep there (the "old" report).
$ws.Rows.Item(95).Insert()

# New row 95 holds the updated weekly report values.
$ws.Range("A95").Value = 5
$ws.Range("B95").Value = "Macroferia Regional de Talca"
$ws.Range("C95").Value = "Maule"
$ws.Range("D95").Value = 44595
$ws.Range("D95").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E95").Value = 7
$ws.Range("F95").Value = 100112001
$ws.Range("G95").Value = "Berenjena"
$ws.Range("H95").Value = "Sin especificar"
$ws.Range("I95").Value = "Primera"
$ws.Range("J95").Value = 150
$ws.Range("K95").Value = 7000
$ws.Range("L95").Value = 7000
$ws.Range("M95").Value = 7000
$ws.Range("N95").Value = "$/caja 50 unidades"
$ws.Range("O95").Value = "Región del Maule"
$ws.Range("P95").Value = 140
$ws.Range("Q95").Value = 50
$ws.Range("R95").Value = "Hortaliza"
